# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that System / system comes first in the comma-separated list, keeping the
# relative order of the remaining entries unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value -split ','
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq 'system') {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -gt 0) {
        $newOrder = $systemParts + $otherParts
        $newValue = [string]::Join(', ', $newOrder)
        if ($newValue -ne $value) {
            $cell.Value = $newValue
        }
    }
}
